# Updated cryptos list with GitHub Actions
# Applies the latest price/volume snapshot to the cryptos worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value  = "67.904.48"
$ws.Range("E2").Value  = "  +1.35%  "

$ws.Range("D3").Value  = "2.628.81"
$ws.Range("E3").Value  = "  +0.82%  "

$ws.Range("E4").Value  = "  -0.13%  "

$ws.Range("D5").Value  = "597.13"
$ws.Range("E5").Value  = "  +0.78%  "

$ws.Range("D6").Value  = "153.74"
$ws.Range("E6").Value  = "  +1.20%  "

$ws.Range("E7").Value  = "  +0.03%  "

$ws.Range("E8").Value  = "  -1.05%  "

$ws.Range("D9").Value  = "2.627.16"
$ws.Range("E9").Value  = "  +0.76%  "

$ws.Range("E10").Value = "  +10.36%  "

$ws.Range("E11").Value = "  -0.57%  "

$ws.Range("E12").Value = "  +1.33%  "

$ws.Range("E13").Value = "  +0.98%  "

$ws.Range("E14").Value = "  +0.59%  "

$ws.Range("E15").Value = "  +5.03%  "

$ws.Range("D16").Value = "3.107.16"
$ws.Range("E16").Value = "  +0.86%  "

$ws.Range("D17").Value = "67.806.73"
$ws.Range("E17").Value = "  +1.39%  "

$ws.Range("D18").Value = "2.625.17"
$ws.Range("E18").Value = "  +0.21%  "

$ws.Range("D19").Value = "11.42"
$ws.Range("E19").Value = "  +4.15%  "

$ws.Range("D20").Value = "370.86"
$ws.Range("E20").Value = "  +1.98%  "

$ws.Range("E21").Value = "  +0.98%  "

$ws.Range("D23").Value = "4.80"
$ws.Range("E23").Value = "  -1.01%  "

$ws.Range("D24").Value = "2.08"
$ws.Range("E24").Value = "  +1.63%  "

$ws.Range("D25").Value = "72.10"
$ws.Range("E25").Value = "  +8.77%  "

$ws.Range("E26").Value = "  +0.18%  "

$ws.Range("D27").Value = "9.88"
$ws.Range("E27").Value = "  -0.67%  "

$ws.Range("E29").Value = "  +3.65%  "

$ws.Range("E30").Value = "  -0.07%  "

$ws.Range("D31").Value = "576.11"
$ws.Range("E31").Value = "  -0.37%  "

$ws.Range("E32").Value = "  +2.45%  "

$ws.Range("E33").Value = "  +1.71%  "

$ws.Range("E34").Value = "  +1.67%  "

$ws.Range("D35").Value = "1.00"

$ws.Range("E36").Value = "  +3.49%  "

$ws.Range("D37").Value = "1.54"
$ws.Range("E37").Value = "  +3.28%  "

$ws.Range("D38").Value = "159.26"
$ws.Range("E38").Value = "  +1.90%  "

$ws.Range("E39").Value = "  +1.09%  "

$ws.Range("E40").Value = "  +5.50%  "

$ws.Range("E41").Value = "  +0.80%  "

$ws.Range("E42").Value = "  +2.38%  "

$ws.Range("E43").Value = "  +4.05%  "

$ws.Range("E44").Value = "  +16.04%  "

$ws.Range("E45").Value = "  +6.16%  "

$ws.Range("E46").Value = "  +0.06%  "

$ws.Range("D47").Value = "40.08"
$ws.Range("E47").Value = "  -2.44%  "

$ws.Range("D48").Value = "155.45"
$ws.Range("E48").Value = "  -0.01%  "

$ws.Range("E49").Value = "  -0.54%  "

$ws.Range("D50").Value = "21.95"
$ws.Range("E50").Value = "  +1.71%  "

$ws.Range("E51").Value = "  +0.34%  "
